# This script applies a cyclic re-shuffle of observation records across
# rows 7-10 and 12-14 of the active worksheet, matching the upstream
# "Automatic update of files." commit.
#
#   Row 7  <-> Row 8   (simple swap)
#   Row 9  <-> Row 10  (simple swap)
#   Row 12 -> Row 13 -> Row 14 -> Row 12  (3-way cycle)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Worksheet, $Row, $Values)
    foreach ($col in $Values.Keys) {
        $Worksheet.Range("$col$Row").Value = $Values[$col]
    }
}

# Snapshot of the "before" values for every row involved in the shuffle,
# keyed by column letter, only for the cells that actually differ between
# rows (the location/date/author columns are identical across the block
# and are intentionally left untouched).
$row7 = @{
    A  = 131042226
    B  = 79244
    E  = 6425
    F  = "Garnlav"
    G  = "Alectoria sarmentosa"
    H  = "(Ach.) Ach."
    Q  = 479114
    R  = 6792438
    S  = 50
    Z  = "16:15"
    AB = "16:15"
    AC = "Rikligt till måttligt i en radie av ca 50 meter,synfältet"
}

$row8 = @{
    A  = 131038653
    B  = 57881
    E  = 100049
    F  = "Spillkråka"
    G  = "Dryocopus martius"
    H  = "(Linnaeus, 1758)"
    M  = "äldre spår"
    Q  = 479094
    R  = 6792753
    S  = 10
    Z  = "12:08"
    AB = "12:08"
}

$row9 = @{
    A = 131039523
    B = 79244
    E = 6425
    F = "Garnlav"
    G = "Alectoria sarmentosa"
    H = "(Ach.) Ach."
    Q = 479079
    R = 6792517
}

$row10 = @{
    A = 131040374
    B = 79002
    E = 228912
    F = "Mörk kolflarnlav"
    G = "Carbonicola myrmecina"
    H = "(Ach.) Bendiksby & Timdal"
    Q = 479088
    R = 6792211
}

$row12 = @{
    A = 131040483
    B = 78647
    D = "NT"
    E = 6437
    F = "Blanksvart spiklav"
    G = "Calicium denigratum"
    H = "(Vain.) Tibell"
    Q = 479088
    R = 6792211
}

$row13 = @{
    A = 131039119
    B = 79244
    D = "NT"
    E = 6425
    F = "Garnlav"
    G = "Alectoria sarmentosa"
    H = "(Ach.) Ach."
    Q = 479105
    R = 6792638
}

$row14 = @{
    A = 131039519
    B = 8451
    D = "LC"
    E = 106545
    F = "Mindre märgborre"
    G = "Tomicus minor"
    H = "(Hartig, 1834)"
    M = "färska gnagspår"
    Q = 479079
    R = 6792517
}

# --- Row 7 <-> Row 8 -------------------------------------------------------
Set-RowValues $ws 7 $row8
$ws.Range("M7").Value = "äldre spår"
$ws.Range("AC7").Value = ""

Set-RowValues $ws 8 $row7
$ws.Range("M8").Value = ""
$ws.Range("AC8").Value = "Rikligt till måttligt i en radie av ca 50 meter,synfältet"

# --- Row 9 <-> Row 10 -------------------------------------------------------
Set-RowValues $ws 9 $row10
Set-RowValues $ws 10 $row9

# --- Row 12 -> Row 13 -> Row 14 -> Row 12 (cycle) --------------------------
Set-RowValues $ws 12 $row14
$ws.Range("M12").Value = "färska gnagspår"

Set-RowValues $ws 13 $row12

Set-RowValues $ws 14 $row13
$ws.Range("M14").Value = ""
